$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.064.64'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").Value = '1.668.68'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = "'216.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").Value = "'0.5108"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = "'0.2645"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").Value = "'0.06408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.85%  '
$ws.Range("D10").Value = "'21.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.29%  '
$ws.Range("D11").Value = "'0.07428"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = '1.682.45'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = "'4.499"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = "'0.5844"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = "'0.000008517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = "'64.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '25.962.17'
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("D18").Value = "'4.940"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = "'10.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").Value = "'190.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("D22").Value = "'6.225"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("D23").Value = "'1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = "'145.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").Value = "'7.600"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = "'0.1198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.61%  '
$ws.Range("D27").Value = "'15.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = "'0.06737"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +18.80%  '
$ws.Range("D29").Value = "'1.319"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("D30").Value = "'1.315"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").Value = "'3.526"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").Value = "'3.519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").Value = "'1.642"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").Value = "'1.019"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D35").Value = "'0.6097"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").Value = "'2.367"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = "'2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("D38").Value = "'6.229"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.18%  '
$ws.Range("D39").Value = "'0.01603"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").Value = '1.076.62'
$ws.Range("E40").Value = '  -2.12%  '
$ws.Range("D41").Value = "'0.8640"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("D42").Value = "'1.009"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").Value = "'100.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.39%  '
$ws.Range("D44").Value = '1.815.32'
$ws.Range("E44").Value = '  -2.32%  '
$ws.Range("D45").Value = "'0.00000000115"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.20%  '
$ws.Range("D46").Value = "'56.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").Value = "'1.009"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("D48").Value = "'8.075"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.98%  '
$ws.Range("D49").Value = "'0.05217"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").Value = "'0.4283"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("D51").Value = "'5.978"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.92%  '
